# Airflow scheduler change: webscrape now captures extra "Unnamed: 0.x" index
# columns on the link sheet, and the prices sheet gains an "updated_time"
# column plus refreshed prices.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "prixworkshop": insert four columns before the existing "link"
# column (old C) so the scraped index columns (Unnamed: 0 .. Unnamed: 0.4)
# all line up, each holding the same running index as column B.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("prixworkshop")

# Shift the old "link" column (C) four slots to the right -> G, freeing up
# C:F for the duplicated index columns.
$ws1.Range("C1:C4").Cut($ws1.Range("G1:G4"))

# Header row.
$ws1.Range("B1").Value = "Unnamed: 0.4"
$ws1.Range("C1").Value = "Unnamed: 0.3"
$ws1.Range("D1").Value = "Unnamed: 0.2"
$ws1.Range("E1").Value = "Unnamed: 0.1"
$ws1.Range("F1").Value = "Unnamed: 0"
$ws1.Range("G1").Value = "link"

for ($r = 2; $r -le 4; $r++) {
    $idx = $ws1.Cells.Item($r, 2).Value2
    $ws1.Cells.Item($r, 3).Value = $idx
    $ws1.Cells.Item($r, 4).Value = $idx
    $ws1.Cells.Item($r, 5).Value = $idx
    $ws1.Cells.Item($r, 6).Value = $idx
}

# Header cells B1:G1 carry the bold/bordered/centered header style already
# used on this sheet - make sure the newly introduced cells pick it up too.
$ws1.Range("B1").Copy()
$ws1.Range("C1:G1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Sheet "prixworkshop_prices": add an "updated_time" column (E) populated by
# the new scheduler run, and refresh the regular-price text in column D.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("prixworkshop_prices")

$ws2.Range("E1").Value = "updated_time"
$ws2.Range("B1").Copy()
$ws2.Range("E1").PasteSpecial(-4122)

$ws2.Range("D2").Value = "`$155.00 CAD"
$ws2.Range("D3").Value = "`$127.00 CAD"
$ws2.Range("D4").Value = "`$20.00 CAD"

$ws2.Range("E2").Value = 45768.43399544809
$ws2.Range("E3").Value = 45768.43402195397
$ws2.Range("E4").Value = 45768.4340465211

# Match the scraper's pandas/openpyxl datetime round-trip: a lower-case
# format gets registered first (numFmtId 164, left unused on any cell) and
# the upper-case variant actually applied to the timestamp cells (165).
$ws2.Range("E2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws2.Range("E2:E4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
